$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.221.42'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '2.254.66'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'307.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.42%  '
$ws.Range('D6').Value = "'98.98"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  -3.19%  '
$ws.Range('D10').Value = "'35.52"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.34%  '
$ws.Range('D11').Value = "'0.0823"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = "'7.32"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.15%  '
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '2.595.92'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').Value = '2.284.90'
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').Value = "'0.840"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '44.082.73'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = "'12.95"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.88%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').Value = "'6.35"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('D22').Value = "'65.46"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').Value = "'243.78"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.16%  '
$ws.Range('E24').Value = '  -6.82%  '
$ws.Range('E25').Value = '  -8.31%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('D29').Value = "'36.49"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = "'6.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = "'20.19"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = "'157.37"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').Value = "'3.57"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.48%  '
$ws.Range('D34').Value = "'0.0827"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.95%  '
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('E37').Value = '  -4.21%  '
$ws.Range('E38').Value = '  -3.63%  '
$ws.Range('D39').Value = "'15.35"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('E40').Value = '  -8.16%  '
$ws.Range('D41').Value = "'3.38"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -10.50%  '
$ws.Range('E42').Value = '  -3.20%  '
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = '1.771.61'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('D45').Value = "'88.63"
$ws.Range('D45').Style = 'Normal'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'16.15"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.92%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').Value = "'5.17"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = "'0.193"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.95%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = "'101.67"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = "'8.28"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = "'70.42"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.58%  '
